# Swap the data (columns B:G) between each pair of adjacent stock-report
# rows. Column A (serial no.) and H:M (blank) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(183,184),
    @(264,265),
    @(313,314),
    @(351,352),
    @(355,356),
    @(382,383),
    @(389,390),
    @(419,420),
    @(421,422),
    @(431,432),
    @(579,580),
    @(583,584),
    @(590,591),
    @(599,600),
    @(601,602),
    @(604,605),
    @(687,688),
    @(709,710),
    @(720,721),
    @(859,860)
)

$cols = @(2,3,4,5,6,7)   # B, C, D, E, F, G

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
